# Auto-generated script to apply cryptos.xlsx diff via Excel COM interop
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $val) {
    $rng = $sheet.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws "D2" "26.259.68"
Set-TextCell $ws "E2" "  +0.30%  "
Set-TextCell $ws "D3" "1.594.85"
Set-TextCell $ws "E3" "  +0.42%  "
Set-TextCell $ws "E4" "  -0.08%  "
Set-TextCell $ws "E5" "  -0.16%  "
Set-TextCell $ws "D6" "0.502"
Set-TextCell $ws "E6" "  +0.01%  "
Set-TextCell $ws "E7" "  -0.06%  "
Set-TextCell $ws "E8" "  +0.23%  "
Set-TextCell $ws "E9" "  +0.03%  "
Set-TextCell $ws "D10" "18.96"
Set-TextCell $ws "E10" "  -1.06%  "
Set-TextCell $ws "E11" "  +1.03%  "
Set-TextCell $ws "D12" "1.819.56"
Set-TextCell $ws "E12" "  +0.42%  "
Set-TextCell $ws "D13" "1.585.01"
Set-TextCell $ws "E13" "  -0.18%  "
Set-TextCell $ws "E14" "  -0.46%  "
Set-TextCell $ws "E15" "  -1.69%  "
Set-TextCell $ws "D16" "63.50"
Set-TextCell $ws "E16" "  -0.20%  "
Set-TextCell $ws "D17" "26.255.27"
Set-TextCell $ws "E17" "  +0.21%  "
Set-TextCell $ws "D18" "230.33"
Set-TextCell $ws "E18" "  +7.75%  "
Set-TextCell $ws "B19" "Chainlink"
Set-TextCell $ws "C19" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell $ws "D19" "7.64"
Set-TextCell $ws "E19" "  +3.77%  "
Set-TextCell $ws "B20" "ShibaInu"
Set-TextCell $ws "C20" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell $ws "D20" "0.0₃0720"
Set-TextCell $ws "E20" "  -0.62%  "
Set-TextCell $ws "D21" "1.00"
Set-TextCell $ws "E21" "  +0.03%  "
Set-TextCell $ws "E22" "  -0.11%  "
Set-TextCell $ws "D23" "2.16"
Set-TextCell $ws "E23" "  +2.13%  "
Set-TextCell $ws "D24" "8.91"
Set-TextCell $ws "E24" "  -0.57%  "
Set-TextCell $ws "D25" "146.13"
Set-TextCell $ws "E25" "  +1.18%  "
Set-TextCell $ws "E26" "  +0.04%  "
Set-TextCell $ws "E27" "  -0.01%  "
Set-TextCell $ws "E28" "  -0.02%  "
Set-TextCell $ws "D29" "15.34"
Set-TextCell $ws "E29" "  +1.81%  "
Set-TextCell $ws "D30" "0.0493"
Set-TextCell $ws "E30" "  -0.26%  "
Set-TextCell $ws "E31" "  -0.07%  "
Set-TextCell $ws "E32" "  +0.60%  "
Set-TextCell $ws "D33" "1.469.55"
Set-TextCell $ws "E33" "  +4.27%  "
Set-TextCell $ws "D34" "2.94"
Set-TextCell $ws "E34" "  -0.62%  "
Set-TextCell $ws "E35" "  -0.36%  "
Set-TextCell $ws "E36" "  +0.43%  "
Set-TextCell $ws "D37" "0.565"
Set-TextCell $ws "E37" "  -3.69%  "
Set-TextCell $ws "E38" "  -0.88%  "
Set-TextCell $ws "D39" "0.817"
Set-TextCell $ws "E39" "  -0.47%  "
Set-TextCell $ws "D40" "5.74"
Set-TextCell $ws "E40" "  -2.60%  "
Set-TextCell $ws "E41" "  -0.03%  "
Set-TextCell $ws "E42" "  +1.60%  "
Set-TextCell $ws "D43" "0.928"
Set-TextCell $ws "E43" "  -1.82%  "
Set-TextCell $ws "D44" "1.731.71"
Set-TextCell $ws "E44" "  +0.48%  "
Set-TextCell $ws "D45" "0.753"
Set-TextCell $ws "E45" "  -1.55%  "
Set-TextCell $ws "D46" "60.35"
Set-TextCell $ws "E46" "  -0.93%  "
Set-TextCell $ws "D47" "88.06"
Set-TextCell $ws "E47" "  +2.57%  "
Set-TextCell $ws "E49" "  +0.02%  "
Set-TextCell $ws "B50" "Algorand"
Set-TextCell $ws "C50" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell $ws "D50" "0.0949"
Set-TextCell $ws "E50" "  -0.77%  "
Set-TextCell $ws "B51" "USDD"
Set-TextCell $ws "C51" "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextCell $ws "D51" "0.999"
Set-TextCell $ws "E51" "  -0.08%  "
